$wb = $excel.ActiveWorkbook

# "Metadata" sheet holds Property/Value pairs (column A = Property, column B = Value)
$ws = $wb.Worksheets.Item("Metadata")

# Row 8: update the "Date" value to the new timestamp
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"

# Row 7: fill in the previously-blank "Experimental" value with the literal text "false".
# A direct Range.Value = "false" assignment gets auto-coerced to a Boolean cell, so
# write it as a formula returning the string, then convert that formula to its
# literal value in place (keeps it plain text and preserves the existing cell style).
$ws.Range("B7").Formula = '="false"'
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)  # xlPasteValues

# Row 17: fill in the previously-blank "Description" value
$ws.Range("B17").Value = "Methods for determining maximum heart rate"
